$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.597878666666666
$ws.Range("H2").Value = 4.793635999999999
$ws.Range("I2").Value = 0.8992131381376172
$ws.Range("J2").Value = 0.8992131381376171
$ws.Range("M2").Value = 0.02035233333333333
$ws.Range("N2").Value = 0.061057
$ws.Range("O2").Value = 0.128389957923202
$ws.Range("P2").Value = 0.128389957923202
$ws.Range("Q2").Value = 0.03252055925022222
$ws.Range("R2").Value = 0.292685033252
$ws.Range("S2").Value = 0.1154499369694791
$ws.Range("T2").Value = 0.115449936969479

$ws.Range("G3").Value = 1.597878666666666
$ws.Range("H3").Value = 4.793635999999999
$ws.Range("I3").Value = 0.8992131381376172
$ws.Range("J3").Value = 0.8992131381376171
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1381673333333333
$ws.Range("N3").Value = 0.414502
$ws.Range("O3").Value = 0.871610042076798
$ws.Range("P3").Value = 0.871610042076798
$ws.Range("Q3").Value = 0.2207746343635555
$ws.Range("R3").Value = 1.986971709272
$ws.Range("S3").Value = 0.7837632011681381
$ws.Range("T3").Value = 0.783763201168138

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1790956666666667
$ws.Range("H4").Value = 0.537287
$ws.Range("I4").Value = 0.1007868618623829
$ws.Range("J4").Value = 0.1007868618623829
$ws.Range("M4").Value = 0.02035233333333333
$ws.Range("N4").Value = 0.061057
$ws.Range("O4").Value = 0.128389957923202
$ws.Range("P4").Value = 0.128389957923202
$ws.Range("Q4").Value = 0.003645014706555555
$ws.Range("R4").Value = 0.032805132359
$ws.Range("S4").Value = 0.01294002095372291
$ws.Range("T4").Value = 0.01294002095372291

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1790956666666667
$ws.Range("H5").Value = 0.537287
$ws.Range("I5").Value = 0.1007868618623829
$ws.Range("J5").Value = 0.1007868618623829
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1381673333333333
$ws.Range("N5").Value = 0.414502
$ws.Range("O5").Value = 0.871610042076798
$ws.Range("P5").Value = 0.871610042076798
$ws.Range("Q5").Value = 0.02474517067488889
$ws.Range("R5").Value = 0.222706536074
$ws.Range("S5").Value = 0.08784684090866003
$ws.Range("T5").Value = 0.08784684090866002
